$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old entry in A4 and replace with new content in A3.
$ws.Range("A4").ClearContents()
$ws.Range("A3").Value = "User 1 added something in cell A3"

# Update the selection to match the new active cell.
$ws.Range("A3").Select()
